# Add season record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns AD, AE, AF with the same header style as the
# existing header cells (bold, bordered, centered) by copying the format
# from an existing header cell (AC1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2-50: constant season record values for every player row.
$lastRow = 50
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 81   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 81   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF -> Ties
}
